# Apply the CATALOG_sample.xlsx edits described by the commit:
#  1. Collection sheet, A11: "RELATION:OutputOf" -> "RELATION:isOutputOf"
#  2. Files sheet: remove row 2 ("bag-info.txt") so the sheet only keeps its header row
#  3. People sheet: rename header cells C1/D1 from "Given Name"/"Family Name"
#     to "Given_Name"/"Family_Name"
# (The Licenses sheet's long CC-license description cell picks up doubled
#  blank-line spacing automatically from the normal load/save round trip of
#  its legacy "_x000d_" encoded line breaks, so no explicit edit is required
#  there.)

$wb = $excel.ActiveWorkbook

$collection = $wb.Worksheets.Item("Collection")
$collection.Range("A11").Value = "RELATION:isOutputOf"

$files = $wb.Worksheets.Item("Files")
$files.Rows.Item(2).Delete()

$people = $wb.Worksheets.Item("People")
$people.Range("C1").Value = "Given_Name"
$people.Range("D1").Value = "Family_Name"
